$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.658923
$ws.Range("H2").Value = 7.976769
$ws.Range("I2").Value = 0.01800502032966059
$ws.Range("J2").Value = 0.01800502032966059
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 24.519512
$ws.Range("N2").Value = 73.558536
$ws.Range("O2").Value = 0.4736537296697991
$ws.Range("P2").Value = 0.4736537296697991
$ws.Range("Q2").Value = 65.19549440557601
$ws.Range("R2").Value = 586.759449650184
$ws.Range("S2").Value = 0.008528145031924294
$ws.Range("T2").Value = 0.008528145031924294

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.658923
$ws.Range("H3").Value = 7.976769
$ws.Range("I3").Value = 0.01800502032966059
$ws.Range("J3").Value = 0.01800502032966059
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 14.70328633333333
$ws.Range("N3").Value = 44.109859
$ws.Range("O3").Value = 0.284029568377475
$ws.Range("P3").Value = 0.284029568377475
$ws.Range("Q3").Value = 39.09490620728567
$ws.Range("R3").Value = 351.854155865571
$ws.Range("S3").Value = 0.00511395815286116
$ws.Range("T3").Value = 0.00511395815286116

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.658923
$ws.Range("H4").Value = 7.976769
$ws.Range("I4").Value = 0.01800502032966059
$ws.Range("J4").Value = 0.01800502032966059
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 12.543947
$ws.Range("N4").Value = 37.631841
$ws.Range("O4").Value = 0.2423167019527259
$ws.Range("P4").Value = 0.2423167019527259
$ws.Range("Q4").Value = 33.35338918908101
$ws.Range("R4").Value = 300.180502701729
$ws.Range("S4").Value = 0.004362917144875136
$ws.Range("T4").Value = 0.004362917144875136

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 122.6832073333333
$ws.Range("H5").Value = 368.049622
$ws.Range("I5").Value = 0.8307550245511554
$ws.Range("J5").Value = 0.8307550245511555
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 24.519512
$ws.Range("N5").Value = 73.558536
$ws.Range("O5").Value = 0.4736537296697991
$ws.Range("P5").Value = 0.4736537296697991
$ws.Range("Q5").Value = 3008.132374408155
$ws.Range("R5").Value = 27073.19136967339
$ws.Range("S5").Value = 0.3934902158205802
$ws.Range("T5").Value = 0.3934902158205803

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 122.6832073333333
$ws.Range("H6").Value = 368.049622
$ws.Range("I6").Value = 0.8307550245511554
$ws.Range("J6").Value = 0.8307550245511555
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 14.70328633333333
$ws.Range("N6").Value = 44.109859
$ws.Range("O6").Value = 0.284029568377475
$ws.Range("P6").Value = 0.284029568377475
$ws.Range("Q6").Value = 1803.8463257137
$ws.Range("R6").Value = 16234.6169314233
$ws.Range("S6").Value = 0.2359589910506833
$ws.Range("T6").Value = 0.2359589910506833

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 122.6832073333333
$ws.Range("H7").Value = 368.049622
$ws.Range("I7").Value = 0.8307550245511554
$ws.Range("J7").Value = 0.8307550245511555
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 12.543947
$ws.Range("N7").Value = 37.631841
$ws.Range("O7").Value = 0.2423167019527259
$ws.Range("P7").Value = 0.2423167019527259
$ws.Range("Q7").Value = 1538.931650579345
$ws.Range("R7").Value = 13850.3848552141
$ws.Range("S7").Value = 0.2013058176798918
$ws.Range("T7").Value = 0.2013058176798918

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 22.334626
$ws.Range("H8").Value = 67.003878
$ws.Range("I8").Value = 0.1512399551191839
$ws.Range("J8").Value = 0.151239955119184
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 24.519512
$ws.Range("N8").Value = 73.558536
$ws.Range("O8").Value = 0.4736537296697991
$ws.Range("P8").Value = 0.4736537296697991
$ws.Range("Q8").Value = 547.6341302225121
$ws.Range("R8").Value = 4928.707172002608
$ws.Range("S8").Value = 0.0716353688172945
$ws.Range("T8").Value = 0.07163536881729451

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 22.334626
$ws.Range("H9").Value = 67.003878
$ws.Range("I9").Value = 0.1512399551191839
$ws.Range("J9").Value = 0.151239955119184
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 14.70328633333333
$ws.Range("N9").Value = 44.109859
$ws.Range("O9").Value = 0.284029568377475
$ws.Range("P9").Value = 0.284029568377475
$ws.Range("Q9").Value = 328.3924012259113
$ws.Range("R9").Value = 2955.531611033202
$ws.Range("S9").Value = 0.0429566191739305
$ws.Range("T9").Value = 0.04295661917393051

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 22.334626
$ws.Range("H10").Value = 67.003878
$ws.Range("I10").Value = 0.1512399551191839
$ws.Range("J10").Value = 0.151239955119184
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 12.543947
$ws.Range("N10").Value = 37.631841
$ws.Range("O10").Value = 0.2423167019527259
$ws.Range("P10").Value = 0.2423167019527259
$ws.Range("Q10").Value = 280.164364808822
$ws.Range("R10").Value = 2521.479283279398
$ws.Range("S10").Value = 0.03664796712795895
$ws.Range("T10").Value = 0.03664796712795895
